# Yamuna.xlsx - "Add files via upload" update
# Refresh the activity log: new dates/classes, a couple of count tweaks,
# and three brand-new rows at the bottom for the latest Apex triggers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 1-9: refresh date / class / count -----------------------

$ws.Range("A1").Value = 43838
$ws.Range("B1").Value = "PortalApplication"
$ws.Range("D1").Value = 2

$ws.Range("A2").Value = 43840
$ws.Range("B2").Value = "CaseTriggerHandler"
$ws.Range("D2").Value = 1

$ws.Range("A3").Value = 43841
$ws.Range("B3").Value = "PortalApplication"
$ws.Range("D3").Value = 1

$ws.Range("A4").Value = 43849
$ws.Range("B4").Value = "RegionUpdateOnCase"
$ws.Range("D4").Value = 1

$ws.Range("A5").Value = 43847
$ws.Range("B5").Value = "Test_EOIUpdate"
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = 43851
$ws.Range("B6").Value = "HttpCreateForms"
$ws.Range("D6").Value = 1

$ws.Range("A7").Value = 43851
$ws.Range("B7").Value = "HttpCreateApprovals"
$ws.Range("D7").Value = 2

$ws.Range("A8").Value = 43853
$ws.Range("B8").Value = "CasePaperTriggerHandler"
$ws.Range("D8").Value = 1

$ws.Range("A9").Value = 43857
$ws.Range("B9").Value = "CashContributionsTriggerHandler"
$ws.Range("D9").Value = 1

# --- Three new rows at the bottom ------------------------------------------

$ws.Range("A10").NumberFormat = "mm/dd/yy;@"
$ws.Range("A10").Value = 43858
$ws.Range("B10").Value = "VolunteerNonCashTriggerHandler"
$ws.Range("C10").Value = "ApexClass"
$ws.Range("D10").Value = 1

$ws.Range("A11").NumberFormat = "mm/dd/yy;@"
$ws.Range("A11").Value = 43858
$ws.Range("B11").Value = "ProjectIncomeTriggerHandler"
$ws.Range("C11").Value = "ApexClass"
$ws.Range("D11").Value = 1

$ws.Range("A12").NumberFormat = "mm/dd/yy;@"
$ws.Range("A12").Value = 43858
$ws.Range("B12").Value = "ProjectCostTriggerHandler"
$ws.Range("C12").Value = "ApexClass"
$ws.Range("D12").Value = 1

# --- Reformat the date column ------------------------------------------

$ws.Range("A1:A12").NumberFormat = "[$-409]d\-mmm\-yyyy;@"

# --- Column A is now manually widened (no longer auto-fit) -----------------

$ws.Columns.Item(1).ColumnWidth = 13.166666666666666

# --- Leave the cursor where the author left it ------------------------------

$ws.Range("H11").Select()
